# Updates crypto price/volume snapshot data (scraped refresh), plus one
# pair of rows (Gas/HuobiToken) that swapped ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''35.127.58'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -0.14%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.903.69'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Value = '''253.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +3.31%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.702'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +2.62%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  -0.42%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''41.51'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +2.54%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  +4.04%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''52.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -1.37%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.0756'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +5.62%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.0978'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -0.15%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''13.17'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +4.54%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''2.181.06'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = '''0.732'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +4.55%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = '''  +5.56%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''1.901.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +0.12%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''35.129.62'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -0.16%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''73.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +2.73%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''0.0₃0842'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +3.59%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''242.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +1.30%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''13.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +3.97%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''5.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +6.25%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = '''  -0.43%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''2.45'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  +6.37%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  -1.32%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''167.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -0.12%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''8.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +1.08%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''18.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +1.94%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.33%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''4.128.77'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -0.33%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '''  +7.07%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''2.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +7.57%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '''  +4.85%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''1.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  +8.49%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''4.24'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +4.20%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  -0.46%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.852'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -6.58%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = '''  +0.26%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''100.53'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +12.43%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''17.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +6.42%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.0215'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +3.75%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''1.11'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +1.86%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''0.0651'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +3.70%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''2.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +0.81%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''1.303.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -3.27%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = '''HuobiToken'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''2.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  +0.08%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = '''Gas'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''12.47'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +2.59%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  -0.94%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  +2.29%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.0752'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +7.21%  '
$ws.Range('E51').Style = 'Normal'
